$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number but must stay text
# (matches the original inlineStr/text storage of the Price column),
# so force text number format before assigning to avoid Excel
# auto-converting them into floating point numbers.
$textCells = @(
    "D5",
    "D6",
    "D8",
    "D10",
    "D12",
    "D13",
    "D14",
    "D16",
    "D19",
    "D20",
    "D21",
    "D23",
    "D25",
    "D26",
    "D28",
    "D29",
    "D30",
    "D31",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D43",
    "D44",
    "D46",
    "D47",
    "D49",
    "D50",
    "D51"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.153.06"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "2.626.80"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "598.73"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "152.23"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.555"
$ws.Range("E8").Value = "  +3.17%  "
$ws.Range("D9").Value = "2.625.98"
$ws.Range("D10").Value = "0.122"
$ws.Range("E10").Value = "  +3.30%  "
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").Value = "5.18"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "0.349"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").Value = "27.52"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").Value = "3.105.89"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").Value = "0.0000181"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").Value = "67.259.73"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "2.629.94"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").Value = "11.15"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").Value = "363.94"
$ws.Range("E20").Value = "  +2.26%  "
$ws.Range("D21").Value = "7.46"
$ws.Range("E21").Value = "  -3.53%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").Value = "2.11"
$ws.Range("E23").Value = "  +4.12%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("B25").Value = "Aptos"
$ws.Range("C25").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D25").Value = "10.12"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "67.51"
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("D27").Value = "2.764.62"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "1.01"
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0000102"
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").Value = "577.88"
$ws.Range("E30").Value = "  -5.88%  "
$ws.Range("D31").Value = "1.39"
$ws.Range("E31").Value = "  -3.31%  "
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("D33").Value = "1.84"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "0.128"
$ws.Range("E35").Value = "  -3.79%  "
$ws.Range("D36").Value = "1.53"
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("D37").Value = "4.91"
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("D38").Value = "157.85"
$ws.Range("E38").Value = "  +2.33%  "
$ws.Range("D39").Value = "19.22"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "0.368"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "5.25"
$ws.Range("E41").Value = "  -3.14%  "
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").Value = "2.55"
$ws.Range("E43").Value = "  +1.68%  "
$ws.Range("D44").Value = "41.24"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "16.35"
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("D47").Value = "155.22"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").Value = "0.0₆0284"
$ws.Range("E48").Value = "  -1.94%  "
$ws.Range("D49").Value = "3.73"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "0.622"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "20.41"
$ws.Range("E51").Value = "  -1.63%  "
